$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(34.913776369470838, 16.046616270646645, 0.45960700729807219, 30.658450542146987, 15.259614157230043, 0.49772946405925655, 338.5, 301)
    3  = @(34.922945321561237, 16.057267664509535, 0.45979133537158634, 30.660296317883805, 15.260272694877289, 0.49772097884051247, 338.5, 301)
    4  = @(34.907910898952814, 16.043983620324116, 0.45960881665953296, 30.655801210163922, 15.266772246831055, 0.49800597747121877, 338.5, 301)
    5  = @(34.905703523018794, 16.045393555512156, 0.45967827420899615, 30.647876950116746, 15.241345824614058, 0.4973051102176263,  338.5, 301)
    6  = @(34.917519062778261, 16.066660234786756, 0.46013178100942631, 30.636079126776572, 15.252227690426082, 0.49785181802508521, 338.5, 301)
    7  = @(34.901535472372629, 16.0713571228364,   0.4604770794556638,  30.965267588605251, 15.541494007393675, 0.50190084626017284, 338.5, 302)
    8  = @(34.887273081492978, 16.104782033631658, 0.46162341195348205, 30.912421485517669, 15.537837079453327, 0.50264056753796316, 338.5, 302)
    9  = @(34.851684462076136, 16.170131056300125, 0.46396985700635662, 30.763331295120999, 15.463066326108835, 0.50264602938373082, 338.5, 302.5)
    10 = @(35.267866443531076, 16.763699684893659, 0.47532503027181278, 31.193345126036419, 16.069243005674345, 0.51514971994015768, 339,   304.5)
    11 = @(35.115137281107216, 17.18384179534123,  0.48935710140556815, 31.782345969065368, 16.801429943956588, 0.52864033260193832, 339,   328)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # Column B is index 2
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
